# REPORTGEN-991: update chinese templates
#
# This script reproduces, via Excel COM automation, the structural/content
# changes made to the "Function-points-sample.xlsx" template:
#
#  - On the "Transactional Functions" and "Data Functions" sheets, the
#    summary table header row (row 6) gains a new "Previous contributed
#    Value" column, and the old "Value" column header is renamed to
#    "Contributed Value".
#  - The RepGen table placeholder tags in both sheets are updated to
#    request the previous-period/zero-suppressed variant of the table.
#  - The "Data Functions" sheet becomes the active / selected sheet
#    instead of "Summary".
#
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Transactional Functions sheet
# ---------------------------------------------------------------------
$wsTF = $wb.Worksheets.Item("Transactional Functions")

# Rename existing "Value" header (column D) to "Contributed Value"
$wsTF.Range("D6").Value = "Contributed Value"

# Insert a new column before the old column E, shifting Description /
# Object Type / Module / Technology one column to the right
$wsTF.Range("E1").EntireColumn.Insert()

# Populate the header of the newly inserted column
$wsTF.Range("E6").Value = "Previous contributed Value"

# Update the RepGen table directive for the transactional functions table
$wsTF.Range("B7").Value = "RepGen:TABLE;IFPUG_FUNCTIONS;HEADER=NO,ZERO=NO,PREVIOUS=YES,TYPE=TF"

# ---------------------------------------------------------------------
# Data Functions sheet
# ---------------------------------------------------------------------
$wsDF = $wb.Worksheets.Item("Data Functions")

# Rename existing "Value" header (column D) to "Contributed Value"
$wsDF.Range("D6").Value = "Contributed Value"

# Insert a new column before the old column E, shifting Description /
# Object Type / Module / Technology one column to the right
$wsDF.Range("E1").EntireColumn.Insert()

# Populate the header of the newly inserted column
$wsDF.Range("E6").Value = "Previous contributed Value"

# Update the RepGen table directive for the data functions table
$wsDF.Range("B7").Value = "RepGen:TABLE;IFPUG_FUNCTIONS;HEADER=NO,ZERO=NO,PREVIOUS=YES,TYPE=DF"

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
# Restore per-sheet selections
$wsTF.Activate()
$wsTF.Range("E6").Select()

$wsDF.Activate()
$wsDF.Range("E6").Select()

# "Data Functions" ends up as the active / selected tab
$wsDF.Activate()
